$wb = $excel.ActiveWorkbook

# --- ev_charging_uc sheet: update the two permutation strings ---
$wsUc = $wb.Worksheets.Item("ev_charging_uc")
$wsUc.Range("C13").Value = "WaP,SaD,FaD,FaP,SaP,RaP,WaD,RaD"
$wsUc.Range("C14").Value = "FaN,SaN,WaN,RaN,FaP,SaP,RaP,WaP"

# --- re_profiles sheet: rotate the M4:N7 ranking table by one row ---
$wsRe = $wb.Worksheets.Item("re_profiles")
$wsRe.Range("M4").Value = "F"
$wsRe.Range("N4").Value = 0.26702915316982878
$wsRe.Range("M5").Value = "S"
$wsRe.Range("N5").Value = 0.4043961129106895
$wsRe.Range("M6").Value = "W"
$wsRe.Range("N6").Value = 0.22555529847292916
$wsRe.Range("M7").Value = "R"
$wsRe.Range("N7").Value = 0.30301943544655252
